# Insert a new row at position 27 (shifts existing rows 27-56 down to 28-57)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eventi")
$ws.Rows("27:27").Insert()

# Populate the new row 27 with the new event record
$ws.Range("A27").Value = 'Mostre'
$ws.Range("B27").Value = 'Modena'
$ws.Range("C27").Value = 'Strada Vaciglio Nord, 6'
$ws.Range("D27").Value = '2022-06-04T08:30:34+00:00'
$ws.Range("F27").Value = '2022-06-04T08:30:59+00:00'
$ws.Range("H27").Value = '2022-06-08T08:00:00+00:00'
$ws.Range("I27").Value = '2022-07-08T09:00:00+00:00'
$ws.Range("J27").Value = 'https://www.comune.modena.it/api/novita/eventi/2022/oro-rosso-fragole-pomodori-molestie-e-sfruttamento-nel-mediterraneo/@@images/04c6eef6-5450-4d2e-a11b-c353a7bdb6b0.jpeg'
$ws.Range("L27").Value = '2022-06-04T08:30:59+00:00'
$ws.Range("M27").Value = 'Sala Renata Bergonzoni della Casa delle Donne'
$ws.Range("N27").Value = ' Inaugurazione mercoledì 8 giugno ore 18.30  mostra aperta dal 10 giugno nei seguenti orari:  venerdì e sabato dalle 10 alle 13 (ad esclusione di venerdì 17 giugno) '
$ws.Range("P27").Value = ' ingresso libero'
$ws.Range("S27").Value = 'Oro rosso. Fragole, pomodori, molestie e sfruttamento nel Mediterraneo'
$ws.Range("X27").Value = 'https://www.comune.modena.it/novita/eventi/2022/oro-rosso-fragole-pomodori-molestie-e-sfruttamento-nel-mediterraneo'
$ws.Range("Y27").Value = '44,64582'
$ws.Range("Z27").Value = '10,92572'
$ws.Range("AA27").Value = 'POINT (10.92572 44.64582)'

# Boolean / numeric fields
$ws.Range("V27").Value = $false
$ws.Range("W27").Value = 41123

# Empty text fields (descrizione, email, desc_img, patrocinato_da, reperibilita, telefono, extrainfo, web)
$ws.Range("E27").Value = ""
$ws.Range("G27").Value = ""
$ws.Range("K27").Value = ""
$ws.Range("O27").Value = ""
$ws.Range("Q27").Value = ""
$ws.Range("R27").Value = ""
$ws.Range("T27").Value = ""
$ws.Range("U27").Value = ""
